$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: column group headers ("Input" x5, "Ouput " for output column) ---
$ws.Range("C5").Value = "Input"
$ws.Range("F5").Value = "Input"
$ws.Range("I5").Value = "Input"
$ws.Range("M5").Value = "Input"
$ws.Range("P5").Value = "Input"
$ws.Range("S5").Value = "Ouput "

# --- Row 7: column titles ---
$ws.Range("A7").Value = "Test Cases"
$ws.Range("C7").Value = "orb_number (0-20)"
$ws.Range("F7").Value = "potion_portion (0.0 -1.0)"
$ws.Range("I7").Value = "word_choice ('Whisper', 'Break', 'Open')"
$ws.Range("M7").Value = "stair_choice ('Left', 'Right')"
$ws.Range("P7").Value = "riddle_choice (1 or 2)"
$ws.Range("S7").Value = "Outcome (Escaped or Failed to escape"

# --- Row 8: test case 1 ---
$ws.Range("A8").Value = 1
$ws.Range("C8").Value = 7
$ws.Range("F8").Value = 0.7
$ws.Range("I8").Value = "Open"
$ws.Range("M8").Value = "Right"
$ws.Range("P8").Value = 1
$ws.Range("S8").Value = "Failed to escape"

# --- Row 9: test case 2 ---
$ws.Range("A9").Value = 2
$ws.Range("C9").Value = 10
$ws.Range("F9").Value = 0.5
$ws.Range("I9").Value = "Break"
$ws.Range("M9").Value = "Right"
$ws.Range("P9").Value = 2
$ws.Range("S9").Value = "Espaced"

# --- Row 10: test case 3 ---
$ws.Range("A10").Value = 3
$ws.Range("C10").Value = 5
$ws.Range("F10").Value = 0.3
$ws.Range("I10").Value = "Whisper"
$ws.Range("M10").Value = "Left"
$ws.Range("P10").Value = 1
$ws.Range("S10").Value = "Espaced"

# --- Row 11: test case 4 ---
$ws.Range("A11").Value = 4
$ws.Range("C11").Value = 15
$ws.Range("F11").Value = 1
$ws.Range("I11").Value = "Open"
$ws.Range("M11").Value = "Left"
$ws.Range("P11").Value = 2
$ws.Range("S11").Value = "Failed to escape"

# --- Move the selection to A7, matching the saved cursor position ---
$ws.Range("A7").Select()
